$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Hamideh%Abbaspour Kasgari%NULL%0, Siavash%Moradi%NULL%1, Amir Mohammad%Shabani%NULL%1, Farhang%Babamahmoodi%NULL%2, Farhang%Babamahmoodi%NULL%0, Ali Reza%Davoudi Badabi%NULL%1, Lotfollah%Davoudi%NULL%1, Ahmad%Alikhani%NULL%1, Akbar%Hedayatizadeh Omran%NULL%1, Majid%Saeedi%NULL%2, Shahin%Merat%NULL%3, Hannah%Wentzel%NULL%6, Hannah%Wentzel%NULL%0, Anna%Garratt%NULL%2, Jacob%Levi%NULL%3, Bryony%Simmons%NULL%3, Andrew%Hill%NULL%4, Hafez%Tirgar Fakheri%hafezfakheri@gmail.com%1]"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "Oxford University Press"

$ws.Range("C3").Value = "Unknown Title"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").Value = "not found"
$ws.Range("G3").Value = "N/A"
$ws.Range("H3").Value = "1970-01-01"
$ws.Range("J3").Value = ""

$ws.Range("E4").Value = "[Gholamali%Eslami%NULL%0, Sajedeh%Mousaviasl%NULL%1, Esmat%Radmanesh%NULL%1, Saeed%Jelvay%NULL%1, Saeid%Bitaraf%NULL%1, Bryony%Simmons%NULL%0, Hannah%Wentzel%NULL%0, Andrew%Hill%NULL%0, Anahita%Sadeghi%NULL%0, James%Freeman%NULL%1, Shokrollah%Salmanzadeh%NULL%1, Hani%Esmaeilian%NULL%1, Morteza%Mobarak%NULL%1, Ramin%Tabibi%NULL%1, Amir Hosein%Jafari Kashi%NULL%1, Zahra%Lotfi%NULL%2, Seyed Mehdi%Talebzadeh%NULL%1, Aseni%Wickramatillake%NULL%1, Mahboobeh%Momtazan%NULL%1, Majid%Hajizadeh Farsani%NULL%1, Sedigheh%Marjani%NULL%1, Sara%Mobarak%s.mobarak@abadanums.ac.ir%1]"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "Oxford University Press"

$ws.Range("E5").Value = "[Fatemeh%Roozbeh%NULL%0, Majid%Saeedi%NULL%0, Reza%Alizadeh-Navaei%NULL%1, Akbar%Hedayatizadeh-Omran%NULL%1, Shahin%Merat%NULL%0, Hannah%Wentzel%NULL%0, Hannah%Wentzel%NULL%0, Jacob%Levi%NULL%0, Andrew%Hill%NULL%0, Amir%Shamshirian%shamshirian.amir@gmail.com%1]"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = "Oxford University Press"

$ws.Range("E6").Value = "[Anahid%Nourian%NULL%0, Hossein%Khalili%NULL%1, Zahra%Ahmadinejad%NULL%1, Hamid%Emadi Kouchak%NULL%1, Sirous%Jafari%NULL%1, Sayed Ali%Dehghan Manshadi%NULL%1, Mehrnaz%Rasolinejad%NULL%1, Abbas%Kebriaeezadeh%NULL%1]"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "Mattioli 1885"

$ws.Range("E7").Value = "[Babak%Sayad%NULL%0, Reza%Khodarahmi%rkhodarahmi@mbrc.ac.ir%2, Reza%Khodarahmi%rkhodarahmi@mbrc.ac.ir%0, Farid%Najafi%NULL%1, Ronak%Miladi%NULL%1, Zeinab%Mohseni Afshar%NULL%1, Feizollah%Mansouri%NULL%1, Zohreh%Rahimi%NULL%1, Maria%Shirvani%NULL%1, Mehdi%Salimi%NULL%1, Siavash%Vaziri%NULL%1, Alireza%Janbakhsh%NULL%1, Fatemeh%Khosravi Shadmani%NULL%1, Arezoo%Bozorgomid%NULL%1, Mohammad Hossein%Zamanian%NULL%1, Mandana%Afsharian%NULL%1]"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "Oxford University Press"

$ws.Range("E8").Value = "[Anahita%Sadeghi%NULL%0, Ali%Ali Asgari%NULL%2, Ali%Ali Asgari%NULL%0, Alireza%Norouzi%NULL%1, Zahedin%Kheiri%NULL%1, Amir%Anushirvani%NULL%1, Mahnaz%Montazeri%NULL%1, Hadiseh%Hosamirudsai%NULL%1, Shirin%Afhami%NULL%1, Elham%Akbarpour%NULL%1, Rasoul%Aliannejad%NULL%2, Rasoul%Aliannejad%NULL%0, Amir Reza%Radmard%NULL%1, Amir H%Davarpanah%NULL%1, Jacob%Levi%NULL%0, Hannah%Wentzel%NULL%0, Ambar%Qavi%NULL%1, Anna%Garratt%NULL%0, Bryony%Simmons%NULL%0, Andrew%Hill%NULL%0, Shahin%Merat%merat@tums.ac.ir%0]"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "Oxford University Press"

$ws.Range("D9").Value = "Only a few treatments are approved for coronavirus disease‐2019 (COVID‐19) infections, with continuous debate about their clinical impact.`n Repurposing antiviral treatments might prove the fastest way to identify effective therapy.`n This trial aimed to evaluate the efficacy and safety of sofosbuvir (SOF) plus daclatasvir (DCV) or ravidasvir (RDV) added to standard care (SOC) for patients with moderate and severe COVID‐19 infection.`n Multicentre parallel randomized controlled open‐label trial.`n One hundred and twenty eligible patients with moderate and severe COVID‐19 infection were randomized to one of the study arms.`n Ten days of treatment with SOF plus DCV or RDV in addition to the standard of care compared to SOC.`n Follow up in 7 days.`n Sum of the counted symptoms at 7 and 10 days, mean change in oxygen saturation level, viral negativity, and rate of intensive care unit (ICU) admission.`n Compared to SOC, the SOF‐DCV group experienced a significantly lower sum of the counted symptoms (fever, headache, generalized aches, or respiratory distress) combined with no evidence of deterioration (ICU admission and mechanical ventilation) on Days 7 and 10 of treatment.`n Oxygen saturation also significantly improved among the SOF‐DCV group compared to SOC starting from Day 4. The study also showed positive trends regarding the efficacy of SOF‐DCV with a lower incidence of mortality.`n On the other hand, adding SOF‐RDV to SOC did not show significant improvements in endpoints.`n The results support the efficacy and safety of SOF‐DCV as an add‐on to SOC for the treatment of moderate to severe COVID‐19 infections.`n"
$ws.Range("E9").Value = "[Sherif%Abbass%sherif533@hotmail.com%0, Ehab%Kamal%NULL%2, Ehab%Kamal%NULL%0, Mohsen%Salama%NULL%2, Mohsen%Salama%NULL%0, Tary%Salman%NULL%1, Alyaa%Sabry%NULL%1, Wael%Abdel‐Razek%NULL%1, Sherine%Helmy%NULL%1, Ahmed%Abdelgwad%NULL%1, Neamt%Sakr%NULL%1, Mohamed%Elgazzar%NULL%1, Mohamed%Einar%NULL%1, Mahmoud%Farouk%NULL%1, Mounir%Saif%NULL%1, Ismail%Shehab%NULL%1, Eman%El‐hosieny%NULL%1, Mai%Mansour%NULL%1, Doaa%Mahdi%NULL%1, El‐Sayed%Tharwa%NULL%1, Mostafa%Salah%NULL%1, Ola%Elrouby%NULL%1, Imam%Waked%NULL%1]"
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = "John Wiley and Sons Inc."
